$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.855.38'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.903.16'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("E4").Value = '  -0.52%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.41%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5018'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.01%  '

$ws.Range("E8").Value = '  -0.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07294'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.93%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9103'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.86'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.43%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07656'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.881.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.482'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.619'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.004'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.54%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008708'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.44%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.894.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.69%  '

$ws.Range("E21").Value = '  -2.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.145'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.14%  '

$ws.Range("E23").Value = '  -0.26%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '154.50'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.38%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.863'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.224'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.01%  '

$ws.Range("E27").Value = '  -0.90%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.927'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08973'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.85%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.216'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.238'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.30%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7693'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.638'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.80%  '

$ws.Range("E35").Value = '  +0.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.557'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.099'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5535'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.015'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.90%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05272'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.54%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.975'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.57%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.541'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.19%  '

$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '111.38'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.62'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4795'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.70%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.003'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.637'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.43'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06078'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9002'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.50%  '
